$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-10T03:01:27.013571+00:00"
$ws.Range("K3").Value = "2025-11-10T03:01:27.013608+00:00"
$ws.Range("K4").Value = "2025-11-10T03:01:27.013629+00:00"
$ws.Range("K5").Value = "2025-11-10T03:01:29.757790+00:00"
$ws.Range("K6").Value = "2025-11-10T03:01:29.757826+00:00"
$ws.Range("K7").Value = "2025-11-10T03:01:29.757848+00:00"
$ws.Range("K8").Value = "2025-11-10T03:01:32.071178+00:00"
$ws.Range("K9").Value = "2025-11-10T03:01:34.749013+00:00"
$ws.Range("K10").Value = "2025-11-10T03:01:34.749044+00:00"
$ws.Range("K11").Value = "2025-11-10T03:01:34.749064+00:00"
$ws.Range("K12").Value = "2025-11-10T03:01:37.039071+00:00"
$ws.Range("K13").Value = "2025-11-10T03:01:37.039103+00:00"
$ws.Range("K14").Value = "2025-11-10T03:01:37.039123+00:00"
$ws.Range("K15").Value = "2025-11-10T03:01:37.039140+00:00"
$ws.Range("K16").Value = "2025-11-10T03:01:42.093260+00:00"
$ws.Range("K17").Value = "2025-11-10T03:01:44.750388+00:00"
$ws.Range("K18").Value = "2025-11-10T03:01:47.538425+00:00"
$ws.Range("K19").Value = "2025-11-10T03:01:47.538456+00:00"
$ws.Range("K20").Value = "2025-11-10T03:01:47.538479+00:00"
$ws.Range("K21").Value = "2025-11-10T03:01:49.820511+00:00"
$ws.Range("K22").Value = "2025-11-10T03:01:52.571095+00:00"
$ws.Range("K23").Value = "2025-11-10T03:01:52.571128+00:00"
$ws.Range("K24").Value = "2025-11-10T03:01:54.886852+00:00"
$ws.Range("K25").Value = "2025-11-10T03:01:54.886888+00:00"
$ws.Range("K26").Value = "2025-11-10T03:01:54.886912+00:00"
$ws.Range("K27").Value = "2025-11-10T03:01:57.677184+00:00"
$ws.Range("K28").Value = "2025-11-10T03:01:57.677216+00:00"
$ws.Range("K29").Value = "2025-11-10T03:01:57.677237+00:00"
$ws.Range("K30").Value = "2025-11-10T03:01:57.677254+00:00"
$ws.Range("K31").Value = "2025-11-10T03:01:57.677271+00:00"
$ws.Range("K32").Value = "2025-11-10T03:02:00.358373+00:00"
$ws.Range("K33").Value = "2025-11-10T03:02:00.358403+00:00"
$ws.Range("K34").Value = "2025-11-10T03:02:03.086739+00:00"
$ws.Range("K35").Value = "2025-11-10T03:02:03.086773+00:00"
$ws.Range("K36").Value = "2025-11-10T03:02:03.086791+00:00"
$ws.Range("K37").Value = "2025-11-10T03:02:05.339741+00:00"
$ws.Range("K38").Value = "2025-11-10T03:02:05.339771+00:00"
$ws.Range("K39").Value = "2025-11-10T03:02:05.339789+00:00"
$ws.Range("K40").Value = "2025-11-10T03:02:07.676227+00:00"
$ws.Range("K41").Value = "2025-11-10T03:02:07.676265+00:00"
$ws.Range("K42").Value = "2025-11-10T03:02:07.676286+00:00"
$ws.Range("K43").Value = "2025-11-10T03:02:07.676304+00:00"
$ws.Range("K44").Value = "2025-11-10T03:02:07.676323+00:00"
$ws.Range("K45").Value = "2025-11-10T03:02:07.676339+00:00"
$ws.Range("K46").Value = "2025-11-10T03:02:10.420265+00:00"
$ws.Range("K47").Value = "2025-11-10T03:02:10.420295+00:00"
$ws.Range("K48").Value = "2025-11-10T03:02:15.509065+00:00"
$ws.Range("K49").Value = "2025-11-10T03:02:15.509096+00:00"
$ws.Range("K50").Value = "2025-11-10T03:02:15.509114+00:00"
$ws.Range("K51").Value = "2025-11-10T03:02:17.895401+00:00"
$ws.Range("K52").Value = "2025-11-10T03:02:17.895429+00:00"
